$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update coin prices and 1h volume percentages

$ws.Range("D2").Value = "26.981.70"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "1.559.05"
$ws.Range("E3").Value = "  +0.54%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").Value = "'207.99"
$ws.Range("E5").Value = "  +0.65%  "
$ws.Range("D6").Value = "'0.489"
$ws.Range("E6").Value = "  +0.53%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "'22.02"
$ws.Range("E8").Value = "  -0.27%  "
$ws.Range("E9").Value = "  +0.81%  "
$ws.Range("D10").Value = "'0.0597"
$ws.Range("E10").Value = "  +1.84%  "
$ws.Range("E11").Value = "  -0.34%  "
$ws.Range("D12").Value = "1.779.13"
$ws.Range("E12").Value = "  +0.39%  "
$ws.Range("D13").Value = "1.507.07"
$ws.Range("E13").Value = "  -2.82%  "
$ws.Range("D14").Value = "'3.73"
$ws.Range("E14").Value = "  -0.23%  "
$ws.Range("E15").Value = "  +0.12%  "
$ws.Range("D16").Value = "26.980.60"
$ws.Range("E16").Value = "  +0.19%  "
$ws.Range("D17").Value = "'61.77"
$ws.Range("E17").Value = "  +0.23%  "
$ws.Range("D18").Value = "0.0₃0705"
$ws.Range("E18").Value = "  +1.32%  "
$ws.Range("D19").Value = "'215.62"
$ws.Range("E19").Value = "  -0.78%  "
$ws.Range("E20").Value = "  +1.53%  "
$ws.Range("E21").Value = "  +0.18%  "
$ws.Range("D22").Value = "'4.14"
$ws.Range("E22").Value = "  +2.23%  "
$ws.Range("D23").Value = "'9.20"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("E24").Value = "  -0.89%  "
$ws.Range("D25").Value = "'152.72"
$ws.Range("E25").Value = "  -0.88%  "
$ws.Range("D26").Value = "'6.59"
$ws.Range("E26").Value = "  -0.29%  "
$ws.Range("E27").Value = "  +0.75%  "
$ws.Range("E28").Value = "  +1.43%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("D30").Value = "'0.0475"
$ws.Range("E30").Value = "  +1.68%  "
$ws.Range("E31").Value = "  +3.23%  "
$ws.Range("E32").Value = "  +0.34%  "
$ws.Range("D33").Value = "'3.18"
$ws.Range("E33").Value = "  +3.50%  "
$ws.Range("D34").Value = "1.422.56"
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("E35").Value = "  +10.36%  "
$ws.Range("E36").Value = "  +1.04%  "
$ws.Range("E37").Value = "  +2.16%  "
$ws.Range("E38").Value = "  +0.89%  "
$ws.Range("D39").Value = "'0.532"
$ws.Range("E39").Value = "  +2.10%  "
$ws.Range("D40").Value = "'5.80"
$ws.Range("E40").Value = "  +1.22%  "
$ws.Range("D41").Value = "'0.809"
$ws.Range("E41").Value = "  +0.18%  "
$ws.Range("E42").Value = "  +0.17%  "
$ws.Range("D45").Value = "'64.61"
$ws.Range("E45").Value = "  +0.49%  "
$ws.Range("E46").Value = "  -1.04%  "
$ws.Range("D47").Value = "1.693.77"
$ws.Range("E47").Value = "  +0.40%  "
$ws.Range("D48").Value = "'86.69"
$ws.Range("E48").Value = "  -1.05%  "
$ws.Range("E49").Value = "  +2.55%  "
$ws.Range("D50").Value = "'0.0518"
$ws.Range("E50").Value = "  -0.17%  "
$ws.Range("D51").Value = "'0.0957"
$ws.Range("E51").Value = "  +0.75%  "

# Rows 43 and 44 swap places (WEMIXToken moves up, MXToken moves down)
# with updated price/volume figures
$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  +0.90%  "

$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").Value = "'2.31"
$ws.Range("E44").Value = "  -0.20%  "
